# Bold every populated cell in the "DisplayName"/"DisplayVersion" table
# (A1:B99): header row + every data row. This mirrors the author's edit,
# which introduced a new bold font (fontId=1) and a matching cellXf
# (applyFont="1") referenced by every cell that already holds data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$firstCol  = $usedRange.Column
$lastCol   = $firstCol + $usedRange.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Only touch cells that actually hold a value, so we don't
        # materialise empty cells (e.g. A3, which has no DisplayName
        # entry) that Excel would otherwise create just to carry a
        # style index.
        if ($cell.Value2 -ne $null) {
            $cell.Font.Bold = $true
        }
    }
}
